$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 18, shifting rows 18:25 down to 19:26
$ws.Rows.Item(18).Insert()

# Fill the new row 18 with the same fixed pattern as the surrounding rows,
# plus the new data point (date 44957, volumen 20, prices 5000).
$ws.Range("A18").Value = 10
$ws.Range("B18").Value = "Vega Modelo de Temuco"
$ws.Range("C18").Value = "La Araucanía"
$ws.Range("D18").Value = 44957
$ws.Range("D18").NumberFormat = $ws.Range("D19").NumberFormat
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = 100112017
$ws.Range("G18").Value = "Ramas de apio"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 20
$ws.Range("K18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("M18").Value = 5000
$ws.Range("N18").Value = "$/paquete"
$ws.Range("O18").Value = "Región de La Araucanía"
$ws.Range("P18").Value = 5000
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = "Hortaliza"
